$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").Value = ""
# Row 106
$ws.Range("H106").Value = 2155.25
$ws.Range("J106").Value = 2000
$ws.Range("L106").Value = 2000
$ws.Range("N106").Value = -3262
# Row 112
$ws.Range("H112").Value = 3549.1667
$ws.Range("J112").Value = 3676.7368
$ws.Range("L112").Value = 11030.2104
$ws.Range("N112").Value = -13246.2104
# Row 113
$ws.Range("H113").Value = 2430.9167
$ws.Range("I113").Value = 4000
$ws.Range("J113").Value = 2288.2727
$ws.Range("K113").Value = 4000
$ws.Range("L113").Value = 2288.2727
$ws.Range("M113").Value = -746
$ws.Range("N113").Value = -8796.2727
# Row 132
$ws.Range("H132").Value = 7357.222
$ws.Range("I132").Value = 1533.0741
$ws.Range("J132").Value = 24829.666
$ws.Range("K132").Value = 4599.2223
$ws.Range("L132").Value = 74488.99800000001
$ws.Range("M132").Value = -2069.2223
$ws.Range("N132").Value = -79548.99800000001
# Row 137
$ws.Range("H137").Value = 3835.52
$ws.Range("I137").Value = 2331.4
$ws.Range("K137").Value = 6994.200000000001
$ws.Range("M137").Value = -4444.200000000001
# Row 138
$ws.Range("H138").Value = 5967.183
$ws.Range("I138").Value = 1131.6154
$ws.Range("J138").Value = 7051.017
$ws.Range("K138").Value = 3394.8462
$ws.Range("L138").Value = 21153.051
$ws.Range("M138").Value = 1745.1538
$ws.Range("N138").Value = -31433.051
# Row 141
$ws.Range("H141").Value = 17880.062
$ws.Range("I141").Value = 15475.615
$ws.Range("K141").Value = 46426.845
$ws.Range("M141").Value = -41246.845

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 16
$ws.Range("H16").Value = 16399.334
$ws.Range("I16").Value = 4599
$ws.Range("K16").Value = 4599
$ws.Range("M16").Value = -4312
# Row 19
$ws.Range("H19").Value = 1249.5
$ws.Range("I19").Value = 1249.5
$ws.Range("K19").Value = 1249.5
$ws.Range("M19").Value = -1020.5
# Row 32
$ws.Range("H32").Value = 15769.918
$ws.Range("I32").Value = 15244.698
$ws.Range("K32").Value = 15244.698
$ws.Range("M32").Value = -14957.698
# Row 43
$ws.Range("H43").Value = 19994
$ws.Range("I43").Value = 19994
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 19994
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -19681
$ws.Range("N43").Value = ""
# Row 45
$ws.Range("H45").Value = 2972.8928
$ws.Range("I45").Value = 1653.2
$ws.Range("K45").Value = 1653.2
$ws.Range("M45").Value = -1276.2
# Row 74
$ws.Range("H74").Value = 5814626.5
$ws.Range("I74").Value = 6250569
$ws.Range("J74").Value = 2056.3333
$ws.Range("K74").Value = 6250569
$ws.Range("L74").Value = 2056.3333
$ws.Range("M74").Value = -6249695
$ws.Range("N74").Value = -3804.3333
# Row 77
$ws.Range("H77").Value = 5814626.5
$ws.Range("I77").Value = 6250569
$ws.Range("J77").Value = 2056.3333
$ws.Range("K77").Value = 31252845
$ws.Range("L77").Value = 10281.6665
$ws.Range("M77").Value = -31248477
$ws.Range("N77").Value = -19017.6665
# Row 122
$ws.Range("H122").Value = 5932.136
$ws.Range("I122").Value = 4529.8823
$ws.Range("K122").Value = 13589.6469
$ws.Range("M122").Value = -11139.6469
# Row 132
$ws.Range("H132").Value = 18907.5
$ws.Range("J132").Value = 10916.286
$ws.Range("L132").Value = 32748.858
$ws.Range("N132").Value = -37808.858

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3080.1
$ws.Range("I20").Value = 3420.7
$ws.Range("J20").Value = 2739.5
$ws.Range("K20").Value = 3420.7
$ws.Range("L20").Value = 2739.5
$ws.Range("M20").Value = -3173.7
$ws.Range("N20").Value = -3233.5
# Row 94
$ws.Range("H94").Value = 1328
$ws.Range("I94").Value = 1267.091
$ws.Range("K94").Value = 1267.091
$ws.Range("M94").Value = -816.0909999999999
# Row 99
$ws.Range("H99").Value = 1965
$ws.Range("I99").Value = 1785.8334
$ws.Range("K99").Value = 1785.8334
$ws.Range("M99").Value = -287.8334
# Row 108
$ws.Range("H108").Value = 85658.336
$ws.Range("J108").Value = 85658.336
$ws.Range("L108").Value = 85658.336
$ws.Range("N108").Value = -93338.336
# Row 134
$ws.Range("H134").Value = 1589.5135
$ws.Range("I134").Value = 1248.2812
$ws.Range("J134").Value = 3773.4
$ws.Range("K134").Value = 3744.8436
$ws.Range("L134").Value = 11320.2
$ws.Range("M134").Value = -1209.8436
$ws.Range("N134").Value = -16390.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1594.4166
$ws.Range("I16").Value = 1401.1428
$ws.Range("J16").Value = 1865
$ws.Range("K16").Value = 1401.1428
$ws.Range("L16").Value = 1865
$ws.Range("M16").Value = -1114.1428
$ws.Range("N16").Value = -2439
# Row 39
$ws.Range("H39").Value = 5500
$ws.Range("I39").Value = 5500
$ws.Range("K39").Value = 5500
$ws.Range("M39").Value = -5109
# Row 49
$ws.Range("H49").Value = 5500
$ws.Range("I49").Value = 5500
$ws.Range("K49").Value = 5500
$ws.Range("M49").Value = -5318
# Row 64
$ws.Range("H64").Value = 39269
$ws.Range("J64").Value = 39269
$ws.Range("L64").Value = 39269
$ws.Range("N64").Value = -39765
# Row 67
$ws.Range("H67").Value = 39269
$ws.Range("J67").Value = 39269
$ws.Range("L67").Value = 39269
$ws.Range("N67").Value = -40985
# Row 105
$ws.Range("H105").Value = 1497.5834
$ws.Range("I105").Value = 1434.5
$ws.Range("K105").Value = 1434.5
$ws.Range("M105").Value = 312.5
# Row 113
$ws.Range("H113").Value = 1594.4166
$ws.Range("I113").Value = 1401.1428
$ws.Range("J113").Value = 1865
$ws.Range("K113").Value = 1401.1428
$ws.Range("L113").Value = 1865
$ws.Range("M113").Value = 768.8571999999999
$ws.Range("N113").Value = -6205
# Row 132
$ws.Range("H132").Value = 12826498
$ws.Range("I132").Value = 19609720
$ws.Range("J132").Value = 13744.444
$ws.Range("K132").Value = 58829160
$ws.Range("L132").Value = 41233.33199999999
$ws.Range("M132").Value = -58826630
$ws.Range("N132").Value = -46293.33199999999
# Row 134
$ws.Range("H134").Value = 4352.5884
$ws.Range("J134").Value = 4260.385
$ws.Range("L134").Value = 12781.155
$ws.Range("N134").Value = -17851.155

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 128
$ws.Range("H128").Value = 118499.25
$ws.Range("I128").Value = 118499.25
$ws.Range("K128").Value = 355497.75
$ws.Range("M128").Value = -350517.75
# Row 131
$ws.Range("H131").Value = 14774518
$ws.Range("J131").Value = 12123123
$ws.Range("L131").Value = 36369369
$ws.Range("N131").Value = -36379449

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7050
$ws.Range("I70").Value = 6853.5264
$ws.Range("J70").Value = 7389.364
$ws.Range("K70").Value = 6853.5264
$ws.Range("L70").Value = 7389.364
$ws.Range("M70").Value = -6583.5264
$ws.Range("N70").Value = -7929.364
# Row 73
$ws.Range("H73").Value = 7050
$ws.Range("I73").Value = 6853.5264
$ws.Range("J73").Value = 7389.364
$ws.Range("K73").Value = 6853.5264
$ws.Range("L73").Value = 7389.364
$ws.Range("M73").Value = -5917.5264
$ws.Range("N73").Value = -9261.364

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 591.0833
$ws.Range("I55").Value = 100.166664
$ws.Range("J55").Value = 1082
$ws.Range("K55").Value = 100.166664
$ws.Range("L55").Value = 1082
$ws.Range("M55").Value = 72.833336
$ws.Range("N55").Value = -1428
# Row 93
$ws.Range("H93").Value = 1311.7241
$ws.Range("I93").Value = 1141
$ws.Range("K93").Value = 1141
$ws.Range("M93").Value = 107
# Row 94
$ws.Range("H94").Value = 40000
$ws.Range("J94").Value = 40000
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41352
# Row 132
$ws.Range("H132").Value = 2672.2
$ws.Range("I132").Value = 2633.6025
$ws.Range("J132").Value = 2809.0454
$ws.Range("K132").Value = 7900.8075
$ws.Range("L132").Value = 8427.136200000001
$ws.Range("M132").Value = -5370.8075
$ws.Range("N132").Value = -13487.1362

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1927.4906
$ws.Range("I132").Value = 778.1111
$ws.Range("J132").Value = 2162.5908
$ws.Range("K132").Value = 2334.3333
$ws.Range("L132").Value = 6487.7724
$ws.Range("M132").Value = 195.6667000000002
$ws.Range("N132").Value = -11547.7724
# Row 136
$ws.Range("H136").Value = 3704.45
$ws.Range("I136").Value = 2255.276
$ws.Range("K136").Value = 6765.828
$ws.Range("M136").Value = -4215.828
